# Update cryptos list (price + volume columns, and one PancakeSwap/Decentraland
# row swap) to match the latest scrape.
# Column layout: A=index(unchanged) B=Coin C=Link D=Price E=Volume(1h)
# D-column values are forced to Text (leading apostrophe, like typing '0.99
# into a cell) and the cell style is reset to "Normal" afterwards so that
# Excel doesn't silently reinterpret values such as "0.9997" or "4.975" as
# numbers, and so no new number-format style gets attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''27.931.63'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +1.47%  '
$ws.Cells.Item(3, 4).Value = '''1.778.99'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +1.71%  '
$ws.Cells.Item(4, 4).Value = '''0.9997'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.25%  '
$ws.Cells.Item(5, 4).Value = '''327.41'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.83%  '
$ws.Cells.Item(6, 4).Value = '''0.9984'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.24%  '
$ws.Cells.Item(7, 4).Value = '''0.4581'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +2.24%  '
$ws.Cells.Item(8, 4).Value = '''0.3590'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -0.29%  '
$ws.Cells.Item(9, 4).Value = '''0.07506'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +0.01%  '
$ws.Cells.Item(10, 4).Value = '''41.97'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -0.01%  '
$ws.Cells.Item(11, 4).Value = '''1.108'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +1.18%  '
$ws.Cells.Item(12, 4).Value = '''0.9986'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -0.25%  '
$ws.Cells.Item(13, 4).Value = '''20.89'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +0.72%  '
$ws.Cells.Item(14, 4).Value = '''6.050'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.38%  '
$ws.Cells.Item(15, 4).Value = '''7.228'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +1.37%  '
$ws.Cells.Item(16, 4).Value = '''1.771.15'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +0.98%  '
$ws.Cells.Item(17, 4).Value = '''93.88'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +0.19%  '
$ws.Cells.Item(18, 4).Value = '''0.00001062'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.07%  '
$ws.Cells.Item(19, 4).Value = '''0.06424'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.54%  '
$ws.Cells.Item(20, 4).Value = '''0.9988'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -0.20%  '
$ws.Cells.Item(21, 4).Value = '''17.13'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +1.83%  '
$ws.Cells.Item(22, 4).Value = '''5.803'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.77%  '
$ws.Cells.Item(23, 4).Value = '''27.955.43'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +1.37%  '
$ws.Cells.Item(24, 5).Value = '  +0.67%  '
$ws.Cells.Item(25, 4).Value = '''2.080'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.79%  '
$ws.Cells.Item(26, 4).Value = '''163.32'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.19%  '
$ws.Cells.Item(27, 4).Value = '''20.31'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -0.90%  '
$ws.Cells.Item(28, 4).Value = '''1.977.03'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +1.25%  '
$ws.Cells.Item(29, 4).Value = '''2.196'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +4.90%  '
$ws.Cells.Item(30, 4).Value = '''126.00'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +0.26%  '
$ws.Cells.Item(31, 4).Value = '''1.107'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +2.21%  '
$ws.Cells.Item(32, 4).Value = '''0.09234'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +1.67%  '
$ws.Cells.Item(33, 4).Value = '''3.666'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -0.07%  '
$ws.Cells.Item(34, 4).Value = '''5.555'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +0.13%  '
$ws.Cells.Item(35, 4).Value = '''11.87'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -0.83%  '
$ws.Cells.Item(36, 4).Value = '''0.02298'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +0.14%  '
$ws.Cells.Item(37, 4).Value = '''0.06151'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +2.24%  '
$ws.Cells.Item(38, 4).Value = '''0.2102'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +0.54%  '
$ws.Cells.Item(39, 4).Value = '''0.6333'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(40, 4).Value = '''4.975'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.57%  '
$ws.Cells.Item(41, 4).Value = '''1.189'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -1.15%  '
$ws.Cells.Item(42, 4).Value = '''1.386'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +0.27%  '
$ws.Cells.Item(43, 4).Value = '''7.868'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.96%  '
$ws.Cells.Item(44, 4).Value = '''13.28'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -0.13%  '
$ws.Cells.Item(45, 2).Value = 'Decentraland'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(45, 4).Value = '''0.5935'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.63%  '
$ws.Cells.Item(46, 2).Value = 'PancakeSwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(46, 4).Value = '''3.742'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +0.67%  '
$ws.Cells.Item(47, 4).Value = '''122.98'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +0.47%  '
$ws.Cells.Item(48, 4).Value = '''1.960'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +0.22%  '
$ws.Cells.Item(49, 4).Value = '''0.06928'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +1.07%  '
$ws.Cells.Item(50, 5).Value = '  -1.05%  '
$ws.Cells.Item(51, 4).Value = '''72.58'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +0.29%  '
